# aggiunta percorso di salvataggio per il file excel
# Adds a new worksheet "query7" at the end of the workbook containing
# the "nome"/"cognome" pair for "Vasco"/"Rossi".

$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet so it lands at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "query7"

# Populate the header row and the single data row.
$newSheet.Range("A1").Value = "nome"
$newSheet.Range("B1").Value = "cognome"
$newSheet.Range("A2").Value = "Vasco"
$newSheet.Range("B2").Value = "Rossi"

# Keep the originally active sheet (query1) selected/active.
$wb.Worksheets.Item(1).Activate()
